# Updates the cryptos list (Price and Volume(1h) columns) as scraped by GitHub Actions.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Cell($addr, $value) {
    $ws.Range($addr).Value = $value
}

# Column D ("Price") values are plain-text numbers in this sheet (not real
# numbers), so force them to stay text by using a leading quote-prefix -
# exactly like typing '94.63 into Excel - instead of letting auto-detection
# turn them into floating point numbers.
function Set-TextCell($addr, $value) {
    $ws.Range($addr).Value = "'" + $value
}

# Row 2 - Bitcoin
Set-TextCell "D2" "43.569.60"
Set-Cell "E2" "  -1.35%  "

# Row 3 - Ethereum
Set-TextCell "D3" "2.280.65"
Set-Cell "E3" "  +1.04%  "

# Row 4 - TetherUSD (Volume only)
Set-Cell "E4" "  +0.02%  "

# Row 5 - Solana
Set-TextCell "D5" "94.63"
Set-Cell "E5" "  -4.49%  "

# Row 6 - BNB
Set-TextCell "D6" "265.94"
Set-Cell "E6" "  -2.45%  "

# Row 7 - XRP (Volume only)
Set-Cell "E7" "  -1.05%  "

# Row 9 - Cardano (Volume only)
Set-Cell "E9" "  -3.85%  "

# Row 10 - Avalanche
Set-TextCell "D10" "44.35"
Set-Cell "E10" "  -7.65%  "

# Row 11 - Dogecoin
Set-TextCell "D11" "0.0932"
Set-Cell "E11" "  -1.58%  "

# Row 12 - Polkadot
Set-TextCell "D12" "7.71"
Set-Cell "E12" "  -5.63%  "

# Row 13 - TRON
Set-TextCell "D13" "0.104"
Set-Cell "E13" "  -0.37%  "

# Row 14 - Wrapped liquid staked Ether 2.0
Set-TextCell "D14" "2.620.15"
Set-Cell "E14" "  +1.27%  "

# Row 15 - Chainlink
Set-TextCell "D15" "15.11"
Set-Cell "E15" "  -2.50%  "

# Row 16 - Polygon (Volume only)
Set-Cell "E16" "  +0.55%  "

# Row 17 - Wrapped Ether
Set-TextCell "D17" "2.285.46"
Set-Cell "E17" "  +1.91%  "

# Row 18 - Wrapped BTC
Set-TextCell "D18" "43.487.39"
Set-Cell "E18" "  -1.48%  "

# Row 19 - Shiba Inu (Volume only)
Set-Cell "E19" "  -0.89%  "

# Row 20 - Uniswap
Set-TextCell "D20" "6.17"
Set-Cell "E20" "  -0.91%  "

# Row 21 - Litecoin
Set-TextCell "D21" "72.22"
Set-Cell "E21" "  +1.91%  "

# Row 22 - Immutable X
Set-TextCell "D22" "2.37"
Set-Cell "E22" "  -0.91%  "

# Row 23 - Bitcoin Cash
Set-TextCell "D23" "233.71"
Set-Cell "E23" "  -0.48%  "

# Row 24 - Internet Computer (DFINITY)
Set-TextCell "D24" "8.92"
Set-Cell "E24" "  -11.33%  "

# Row 25 - Dai (Volume only)
Set-Cell "E25" "  -0.05%  "

# Row 26 - PancakeSwap (Volume only)
Set-Cell "E26" "  -0.89%  "

# Row 27 - Cosmos
Set-TextCell "D27" "11.14"
Set-Cell "E27" "  -2.14%  "

# Row 28 - WEMIX Token
Set-TextCell "D28" "3.48"
Set-Cell "E28" "  -1.48%  "

# Row 29 - Toncoin (Price only)
Set-TextCell "D29" "2.28"

# Row 30 - Injective Protocol
Set-TextCell "D30" "39.19"
Set-Cell "E30" "  -2.30%  "

# Row 31 - Monero
Set-TextCell "D31" "175.16"
Set-Cell "E31" "  +0.90%  "

# Row 32 - Ethereum Classic
Set-TextCell "D32" "21.81"
Set-Cell "E32" "  +2.68%  "

# Row 33 - Hedera (Volume only)
Set-Cell "E33" "  -4.12%  "

# Row 34 - Filecoin
Set-TextCell "D34" "5.31"
Set-Cell "E34" "  -7.36%  "

# Row 35 - Stellar (Volume only)
Set-Cell "E35" "  -0.16%  "

# Row 36 - Kaspa (Volume only)
Set-Cell "E36" "  -5.53%  "

# Row 37 - VeChain (Volume only)
Set-Cell "E37" "  -1.30%  "

# Row 38 - RenderToken
Set-TextCell "D38" "4.38"
Set-Cell "E38" "  +0.44%  "

# Row 39 - NEAR Protocol
Set-TextCell "D39" "3.31"
Set-Cell "E39" "  -6.83%  "

# Row 40 - Lido DAO Token
Set-TextCell "D40" "2.33"
Set-Cell "E40" "  +6.27%  "

# Row 41 - Algorand
Set-TextCell "D41" "0.234"
Set-Cell "E41" "  -6.04%  "

# Row 42 - ARBITRUM (Volume only)
Set-Cell "E42" "  +14.79%  "

# Row 43 - Celestia
Set-TextCell "D43" "11.89"
Set-Cell "E43" "  -5.37%  "

# Row 44 - MultiversX
Set-TextCell "D44" "63.16"
Set-Cell "E44" "  +1.59%  "

# Row 45 - FraxShare
Set-TextCell "D45" "8.77"
Set-Cell "E45" "  +2.74%  "

# Row 46 - THORChain
Set-TextCell "D46" "5.20"
Set-Cell "E46" "  -4.75%  "

# Row 47 - Cronos (Volume only)
Set-Cell "E47" "  -1.62%  "

# Row 48 - Aave
Set-TextCell "D48" "97.36"
Set-Cell "E48" "  -3.18%  "

# Row 49 - TrustWalletToken (Volume only)
Set-Cell "E49" "  -1.14%  "

# Row 50 - Stacks
Set-TextCell "D50" "1.49"
Set-Cell "E50" "  +4.47%  "

# Row 51 - RocketPoolETH
Set-TextCell "D51" "2.499.63"
Set-Cell "E51" "  +1.14%  "
